$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text while we
# write the new values, so numeric-looking strings (e.g. "409.49",
# "0.0000220", "61.984.85") are not silently coerced into floating
# point numbers by the COM layer. We restore the default "Normal"
# style afterwards so the cells end up with no explicit style, just
# like in the original workbook.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "62.024.73"
$ws.Range("E2").Value = "  -1.07%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "3.413.41"
$ws.Range("E3").Value = "  -0.78%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  -0.07%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "409.49"
$ws.Range("E5").Value = "  +0.44%  "

# --- Row 6 (Solana) ---
$ws.Range("D6").Value = "129.64"
$ws.Range("E6").Value = "  -0.32%  "

# --- Row 7 (XRP) ---
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +6.48%  "

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = "  +0.02%  "

# --- Row 9 (Cardano) ---
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  +5.88%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  +2.76%  "

# --- Row 11 (Avalanche) ---
$ws.Range("D11").Value = "42.85"
$ws.Range("E11").Value = "  +1.95%  "

# --- Row 12 / Row 13 swap (Polkadot <-> ShibaInu) ---
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "0.0000220"
$ws.Range("E12").Value = "  +42.22%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "9.31"
$ws.Range("E13").Value = "  +10.20%  "

# --- Row 14 (TRON) ---
$ws.Range("E14").Value = "  -0.28%  "

# --- Row 15 (WrappedliquidstakedEther2.0) ---
$ws.Range("D15").Value = "3.954.45"
$ws.Range("E15").Value = "  -0.91%  "

# --- Row 16 (Chainlink) ---
$ws.Range("D16").Value = "21.25"
$ws.Range("E16").Value = "  +6.95%  "

# --- Row 17 (WrappedEther) ---
$ws.Range("D17").Value = "3.428.95"
$ws.Range("E17").Value = "  -0.65%  "

# --- Row 18 (Uniswap) ---
$ws.Range("D18").Value = "12.50"
$ws.Range("E18").Value = "  +8.22%  "

# --- Row 19 (Polygon) ---
$ws.Range("E19").Value = "  +7.06%  "

# --- Row 20 (WrappedBTC) ---
$ws.Range("D20").Value = "62.018.59"
$ws.Range("E20").Value = "  -0.98%  "

# --- Row 21 (BitcoinCash) ---
$ws.Range("D21").Value = "456.27"
$ws.Range("E21").Value = "  +44.89%  "

# --- Row 22 (Litecoin) ---
$ws.Range("D22").Value = "92.09"
$ws.Range("E22").Value = "  +9.01%  "

# --- Row 23 (ImmutableX) ---
$ws.Range("E23").Value = "  +1.41%  "

# --- Row 24 (InternetComputer(DFINITY)) ---
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").Value = "  +3.19%  "

# --- Row 26 (Filecoin) ---
$ws.Range("D26").Value = "9.39"
$ws.Range("E26").Value = "  +15.11%  "

# --- Row 27 (EthereumClassic) ---
$ws.Range("D27").Value = "33.20"
$ws.Range("E27").Value = "  +11.63%  "

# --- Row 28 (LEO) ---
$ws.Range("E28").Value = "  +0.49%  "

# --- Row 29 (RenderToken) ---
$ws.Range("D29").Value = "7.68"
$ws.Range("E29").Value = "  -1.95%  "

# --- Row 30 (Toncoin) ---
$ws.Range("D30").Value = "2.77"
$ws.Range("E30").Value = "  -0.93%  "

# --- Row 31 (Cosmos) ---
$ws.Range("D31").Value = "12.02"
$ws.Range("E31").Value = "  +5.59%  "

# --- Row 32 (Kaspa) ---
$ws.Range("D32").Value = "0.172"
$ws.Range("E32").Value = "  -1.37%  "

# --- Row 33 (Hedera) ---
$ws.Range("E33").Value = "  -0.32%  "

# --- Row 34 (InjectiveProtocol) ---
$ws.Range("D34").Value = "42.95"
$ws.Range("E34").Value = "  -3.62%  "

# --- Row 35 (Dai) ---
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.08%  "

# --- Row 36 (VeChain) ---
$ws.Range("D36").Value = "0.0505"
$ws.Range("E36").Value = "  +4.16%  "

# --- Row 37 (OKB) ---
$ws.Range("D37").Value = "53.82"
$ws.Range("E37").Value = "  +3.97%  "

# --- Row 38 (FirstDigitalUSD) ---
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.10%  "

# --- Row 39 (LidoDAOToken) ---
$ws.Range("E39").Value = "  +1.52%  "

# --- Row 40 (Stellar) ---
$ws.Range("E40").Value = "  +7.84%  "

# --- Row 41 (Stacks) ---
$ws.Range("E41").Value = "  -0.53%  "

# --- Row 42 (TheGraph) ---
$ws.Range("E42").Value = "  -1.73%  "

# --- Row 43 (Monero) ---
$ws.Range("D43").Value = "142.66"
$ws.Range("E43").Value = "  +0.01%  "

# --- Row 44 (NEARProtocol) ---
$ws.Range("D44").Value = "4.26"
$ws.Range("E44").Value = "  +8.62%  "

# --- Row 45 (WEMIXToken) ---
$ws.Range("E45").Value = "  +15.99%  "

# --- Row 47 (Celestia) ---
$ws.Range("D47").Value = "16.64"
$ws.Range("E47").Value = "  -1.32%  "

# --- Row 48 (Cronos) ---
$ws.Range("D48").Value = "0.148"
$ws.Range("E48").Value = "  +22.86%  "

# --- Row 49 (EnergySwap) ---
$ws.Range("D49").Value = "22.49"
$ws.Range("E49").Value = "  +5.71%  "

# --- Row 50 (ThetaToken) ---
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  +9.37%  "

# --- Row 51 (RocketPoolETH) ---
$ws.Range("D51").Value = "3.759.64"
$ws.Range("E51").Value = "  -0.83%  "

# Restore default styling on the touched range (removes the temporary
# text NumberFormat so cells keep their original, implicit style).
$ws.Range("D2:E51").Style = "Normal"
